# Revert "Revert "sms template""
# This re-applies the "sms template" process-variable rows that a prior
# revert had removed: it adds 5 new process variables (sorted into the
# existing VAR_UID-ordered table) and rewrites rows 6-30 accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

    # Row 6: VAR_UID=310614396652148af08f122080329400
    $ws.Cells.Item(6, 1).Value = 28
    $ws.Cells.Item(6, 2).Value = '310614396652148af08f122080329400'
    $ws.Cells.Item(6, 3).Value = '61815949564e4999e9835d8060262100'
    $ws.Cells.Item(6, 4).Value = 2
    $ws.Cells.Item(6, 5).Value = 'agency_manager_national_id'
    $ws.Cells.Item(6, 6).Value = 'string'
    $ws.Cells.Item(6, 7).Value = 1
    $ws.Cells.Item(6, 8).Value = 10
    $ws.Cells.Item(6, 9).Value = 'string'
    $ws.Cells.Item(6, 10).Value = 'workflow'
    $ws.Cells.Item(6, 12).Value = 0
    $ws.Cells.Item(6, 14).Value = '[]'
    # Row 7: VAR_UID=32406339064e34dce406b64089642611
    $ws.Cells.Item(7, 1).Value = 6
    $ws.Cells.Item(7, 2).Value = '32406339064e34dce406b64089642611'
    $ws.Cells.Item(7, 3).Value = '20109551764e348a7a8c913045934777'
    $ws.Cells.Item(7, 4).Value = 1
    $ws.Cells.Item(7, 5).Value = 'substitute'
    $ws.Cells.Item(7, 6).Value = 'string'
    $ws.Cells.Item(7, 7).Value = 1
    $ws.Cells.Item(7, 8).Value = 10
    $ws.Cells.Item(7, 9).Value = 'string'
    $ws.Cells.Item(7, 10).Value = 'workflow'
    $ws.Cells.Item(7, 12).Value = 0
    $ws.Cells.Item(7, 14).Value = '[]'
    # Row 8: VAR_UID=38980047265214996a96853053294934
    $ws.Cells.Item(8, 1).Value = 29
    $ws.Cells.Item(8, 2).Value = '38980047265214996a96853053294934'
    $ws.Cells.Item(8, 3).Value = '61815949564e4999e9835d8060262100'
    $ws.Cells.Item(8, 4).Value = 2
    $ws.Cells.Item(8, 5).Value = 'inactivity_commitment_image'
    $ws.Cells.Item(8, 6).Value = 'file'
    $ws.Cells.Item(8, 7).Value = 8
    $ws.Cells.Item(8, 8).Value = 10
    $ws.Cells.Item(8, 9).Value = 'file'
    $ws.Cells.Item(8, 10).Value = 'workflow'
    $ws.Cells.Item(8, 12).Value = 0
    $ws.Cells.Item(8, 14).Value = '[]'
    $ws.Cells.Item(8, 15).Value = '97290740465214979a6b891095846179'
    # Row 9: VAR_UID=477146009652119e56206b8097951634
    $ws.Cells.Item(9, 1).Value = 21
    $ws.Cells.Item(9, 2).Value = '477146009652119e56206b8097951634'
    $ws.Cells.Item(9, 3).Value = '61815949564e4999e9835d8060262100'
    $ws.Cells.Item(9, 4).Value = 2
    $ws.Cells.Item(9, 5).Value = 'warning_number'
    $ws.Cells.Item(9, 6).Value = 'string'
    $ws.Cells.Item(9, 7).Value = 1
    $ws.Cells.Item(9, 8).Value = 10
    $ws.Cells.Item(9, 9).Value = 'string'
    $ws.Cells.Item(9, 10).Value = 'workflow'
    $ws.Cells.Item(9, 12).Value = 0
    $ws.Cells.Item(9, 14).Value = '[]'
    # Row 10: VAR_UID=47734422964e34b275ab385060972873
    $ws.Cells.Item(10, 1).Value = 1
    $ws.Cells.Item(10, 2).Value = '47734422964e34b275ab385060972873'
    $ws.Cells.Item(10, 3).Value = '20109551764e348a7a8c913045934777'
    $ws.Cells.Item(10, 4).Value = 1
    $ws.Cells.Item(10, 5).Value = 'user_id'
    $ws.Cells.Item(10, 6).Value = 'string'
    $ws.Cells.Item(10, 7).Value = 1
    $ws.Cells.Item(10, 8).Value = 10
    $ws.Cells.Item(10, 9).Value = 'string'
    $ws.Cells.Item(10, 10).Value = 'workflow'
    $ws.Cells.Item(10, 12).Value = 0
    $ws.Cells.Item(10, 14).Value = '[]'
    # Row 11: VAR_UID=4929670796509713bd54bb4006211907
    $ws.Cells.Item(11, 1).Value = 18
    $ws.Cells.Item(11, 2).Value = '4929670796509713bd54bb4006211907'
    $ws.Cells.Item(11, 3).Value = '20109551764e348a7a8c913045934777'
    $ws.Cells.Item(11, 4).Value = 1
    $ws.Cells.Item(11, 5).Value = 'bossResult'
    $ws.Cells.Item(11, 6).Value = 'boolean'
    $ws.Cells.Item(11, 7).Value = 4
    $ws.Cells.Item(11, 8).Value = 10
    $ws.Cells.Item(11, 9).Value = 'boolean'
    $ws.Cells.Item(11, 10).Value = 'workflow'
    $ws.Cells.Item(11, 12).Value = 0
    $ws.Cells.Item(11, 14).Value = '[{"value":"1","label":"موافق"},{"value":"0","label":"مخالف"}]'
    # Row 12: VAR_UID=49694249764f7504d963f32006601246
    $ws.Cells.Item(12, 1).Value = 16
    $ws.Cells.Item(12, 2).Value = '49694249764f7504d963f32006601246'
    $ws.Cells.Item(12, 3).Value = '61815949564e4999e9835d8060262100'
    $ws.Cells.Item(12, 4).Value = 2
    $ws.Cells.Item(12, 5).Value = 'description'
    $ws.Cells.Item(12, 6).Value = 'string'
    $ws.Cells.Item(12, 7).Value = 1
    $ws.Cells.Item(12, 8).Value = 10
    $ws.Cells.Item(12, 9).Value = 'string'
    $ws.Cells.Item(12, 10).Value = 'workflow'
    $ws.Cells.Item(12, 12).Value = 0
    $ws.Cells.Item(12, 14).Value = '[]'
    # Row 13: VAR_UID=51155704065211a50df4688020707686
    $ws.Cells.Item(13, 1).Value = 22
    $ws.Cells.Item(13, 2).Value = '51155704065211a50df4688020707686'
    $ws.Cells.Item(13, 3).Value = '61815949564e4999e9835d8060262100'
    $ws.Cells.Item(13, 4).Value = 2
    $ws.Cells.Item(13, 5).Value = 'warning_registrar'
    $ws.Cells.Item(13, 6).Value = 'string'
    $ws.Cells.Item(13, 7).Value = 1
    $ws.Cells.Item(13, 8).Value = 10
    $ws.Cells.Item(13, 9).Value = 'string'
    $ws.Cells.Item(13, 10).Value = 'workflow'
    $ws.Cells.Item(13, 12).Value = 0
    $ws.Cells.Item(13, 14).Value = '[]'
    # Row 14: VAR_UID=54462855165097103c15870075516338
    $ws.Cells.Item(14, 1).Value = 17
    $ws.Cells.Item(14, 2).Value = '54462855165097103c15870075516338'
    $ws.Cells.Item(14, 3).Value = '20109551764e348a7a8c913045934777'
    $ws.Cells.Item(14, 4).Value = 1
    $ws.Cells.Item(14, 5).Value = 'bossName'
    $ws.Cells.Item(14, 6).Value = 'string'
    $ws.Cells.Item(14, 7).Value = 1
    $ws.Cells.Item(14, 8).Value = 10
    $ws.Cells.Item(14, 9).Value = 'string'
    $ws.Cells.Item(14, 10).Value = 'workflow'
    $ws.Cells.Item(14, 12).Value = 0
    $ws.Cells.Item(14, 14).Value = '[]'
    # Row 15: VAR_UID=55987158364e34bf598e893077989680
    $ws.Cells.Item(15, 1).Value = 2
    $ws.Cells.Item(15, 2).Value = '55987158364e34bf598e893077989680'
    $ws.Cells.Item(15, 3).Value = '20109551764e348a7a8c913045934777'
    $ws.Cells.Item(15, 4).Value = 1
    $ws.Cells.Item(15, 5).Value = 'type'
    $ws.Cells.Item(15, 6).Value = 'string'
    $ws.Cells.Item(15, 7).Value = 1
    $ws.Cells.Item(15, 8).Value = 10
    $ws.Cells.Item(15, 9).Value = 'string'
    $ws.Cells.Item(15, 10).Value = 'workflow'
    $ws.Cells.Item(15, 12).Value = 0
    $ws.Cells.Item(15, 14).Value = '[{"value":"daily","label":"روزانه"},{"value":"hourly","label":"ساعتی"}]'
    # Row 16: VAR_UID=57981362264e35ecfcd2552001673570
    $ws.Cells.Item(16, 1).Value = 11
    $ws.Cells.Item(16, 2).Value = '57981362264e35ecfcd2552001673570'
    $ws.Cells.Item(16, 3).Value = '20109551764e348a7a8c913045934777'
    $ws.Cells.Item(16, 4).Value = 1
    $ws.Cells.Item(16, 5).Value = 'departmentManagerName'
    $ws.Cells.Item(16, 6).Value = 'string'
    $ws.Cells.Item(16, 7).Value = 1
    $ws.Cells.Item(16, 8).Value = 10
    $ws.Cells.Item(16, 9).Value = 'string'
    $ws.Cells.Item(16, 10).Value = 'workflow'
    $ws.Cells.Item(16, 12).Value = 0
    $ws.Cells.Item(16, 14).Value = '[]'
    # Row 17: VAR_UID=58255037264f74fedcb5599050721625
    $ws.Cells.Item(17, 1).Value = 15
    $ws.Cells.Item(17, 2).Value = '58255037264f74fedcb5599050721625'
    $ws.Cells.Item(17, 3).Value = '61815949564e4999e9835d8060262100'
    $ws.Cells.Item(17, 4).Value = 2
    $ws.Cells.Item(17, 5).Value = 'agency_address'
    $ws.Cells.Item(17, 6).Value = 'string'
    $ws.Cells.Item(17, 7).Value = 1
    $ws.Cells.Item(17, 8).Value = 10
    $ws.Cells.Item(17, 9).Value = 'string'
    $ws.Cells.Item(17, 10).Value = 'workflow'
    $ws.Cells.Item(17, 12).Value = 0
    $ws.Cells.Item(17, 14).Value = '[]'
    # Row 18: VAR_UID=64979433964f74fe1f0cd56030135781
    $ws.Cells.Item(18, 1).Value = 14
    $ws.Cells.Item(18, 2).Value = '64979433964f74fe1f0cd56030135781'
    $ws.Cells.Item(18, 3).Value = '61815949564e4999e9835d8060262100'
    $ws.Cells.Item(18, 4).Value = 2
    $ws.Cells.Item(18, 5).Value = 'agency_mobile'
    $ws.Cells.Item(18, 6).Value = 'string'
    $ws.Cells.Item(18, 7).Value = 1
    $ws.Cells.Item(18, 8).Value = 10
    $ws.Cells.Item(18, 9).Value = 'string'
    $ws.Cells.Item(18, 10).Value = 'workflow'
    $ws.Cells.Item(18, 12).Value = 0
    $ws.Cells.Item(18, 14).Value = '[]'
    # Row 19: VAR_UID=68917078964e34cd40bb933062498466
    $ws.Cells.Item(19, 1).Value = 3
    $ws.Cells.Item(19, 2).Value = '68917078964e34cd40bb933062498466'
    $ws.Cells.Item(19, 3).Value = '20109551764e348a7a8c913045934777'
    $ws.Cells.Item(19, 4).Value = 1
    $ws.Cells.Item(19, 5).Value = 'startDate'
    $ws.Cells.Item(19, 6).Value = 'datetime'
    $ws.Cells.Item(19, 7).Value = 5
    $ws.Cells.Item(19, 8).Value = 10
    $ws.Cells.Item(19, 9).Value = 'datetime'
    $ws.Cells.Item(19, 10).Value = 'workflow'
    $ws.Cells.Item(19, 12).Value = 0
    $ws.Cells.Item(19, 14).Value = '[]'
    # Row 20: VAR_UID=81334237864e34de9467553061447577
    $ws.Cells.Item(20, 1).Value = 7
    $ws.Cells.Item(20, 2).Value = '81334237864e34de9467553061447577'
    $ws.Cells.Item(20, 3).Value = '20109551764e348a7a8c913045934777'
    $ws.Cells.Item(20, 4).Value = 1
    $ws.Cells.Item(20, 5).Value = 'emergencyPhone'
    $ws.Cells.Item(20, 6).Value = 'string'
    $ws.Cells.Item(20, 7).Value = 1
    $ws.Cells.Item(20, 8).Value = 10
    $ws.Cells.Item(20, 9).Value = 'string'
    $ws.Cells.Item(20, 10).Value = 'workflow'
    $ws.Cells.Item(20, 12).Value = 0
    $ws.Cells.Item(20, 14).Value = '[]'
    # Row 21: VAR_UID=82864965964e3561d0a6788018415658
    $ws.Cells.Item(21, 1).Value = 10
    $ws.Cells.Item(21, 2).Value = '82864965964e3561d0a6788018415658'
    $ws.Cells.Item(21, 3).Value = '20109551764e348a7a8c913045934777'
    $ws.Cells.Item(21, 4).Value = 1
    $ws.Cells.Item(21, 5).Value = 'departmentManagerResult'
    $ws.Cells.Item(21, 6).Value = 'boolean'
    $ws.Cells.Item(21, 7).Value = 4
    $ws.Cells.Item(21, 8).Value = 10
    $ws.Cells.Item(21, 9).Value = 'boolean'
    $ws.Cells.Item(21, 10).Value = 'workflow'
    $ws.Cells.Item(21, 12).Value = 0
    $ws.Cells.Item(21, 14).Value = '[{"value":"1","label":"موافق"},{"value":"0","label":"مخالف"}]'
    # Row 22: VAR_UID=865431231652119ac176a93032535681
    $ws.Cells.Item(22, 1).Value = 20
    $ws.Cells.Item(22, 2).Value = '865431231652119ac176a93032535681'
    $ws.Cells.Item(22, 3).Value = '61815949564e4999e9835d8060262100'
    $ws.Cells.Item(22, 4).Value = 2
    $ws.Cells.Item(22, 5).Value = 'warning_register_date'
    $ws.Cells.Item(22, 6).Value = 'string'
    $ws.Cells.Item(22, 7).Value = 1
    $ws.Cells.Item(22, 8).Value = 10
    $ws.Cells.Item(22, 9).Value = 'string'
    $ws.Cells.Item(22, 10).Value = 'workflow'
    $ws.Cells.Item(22, 12).Value = 0
    $ws.Cells.Item(22, 14).Value = '[]'
    # Row 23: VAR_UID=86648548864e34d20c5a7c7083843416
    $ws.Cells.Item(23, 1).Value = 4
    $ws.Cells.Item(23, 2).Value = '86648548864e34d20c5a7c7083843416'
    $ws.Cells.Item(23, 3).Value = '20109551764e348a7a8c913045934777'
    $ws.Cells.Item(23, 4).Value = 1
    $ws.Cells.Item(23, 5).Value = 'endDate'
    $ws.Cells.Item(23, 6).Value = 'datetime'
    $ws.Cells.Item(23, 7).Value = 5
    $ws.Cells.Item(23, 8).Value = 10
    $ws.Cells.Item(23, 9).Value = 'datetime'
    $ws.Cells.Item(23, 10).Value = 'workflow'
    $ws.Cells.Item(23, 12).Value = 0
    $ws.Cells.Item(23, 14).Value = '[]'
    # Row 24: VAR_UID=90755691665211a8fac1d72057498881
    $ws.Cells.Item(24, 1).Value = 24
    $ws.Cells.Item(24, 2).Value = '90755691665211a8fac1d72057498881'
    $ws.Cells.Item(24, 3).Value = '61815949564e4999e9835d8060262100'
    $ws.Cells.Item(24, 4).Value = 2
    $ws.Cells.Item(24, 5).Value = 'warning_receiver'
    $ws.Cells.Item(24, 6).Value = 'string'
    $ws.Cells.Item(24, 7).Value = 1
    $ws.Cells.Item(24, 8).Value = 10
    $ws.Cells.Item(24, 9).Value = 'string'
    $ws.Cells.Item(24, 10).Value = 'workflow'
    $ws.Cells.Item(24, 12).Value = 0
    $ws.Cells.Item(24, 14).Value = '[]'
    # Row 25: VAR_UID=92013649665214880f338a7062586547
    $ws.Cells.Item(25, 1).Value = 27
    $ws.Cells.Item(25, 2).Value = '92013649665214880f338a7062586547'
    $ws.Cells.Item(25, 3).Value = '61815949564e4999e9835d8060262100'
    $ws.Cells.Item(25, 4).Value = 2
    $ws.Cells.Item(25, 5).Value = 'referral'
    $ws.Cells.Item(25, 6).Value = 'boolean'
    $ws.Cells.Item(25, 7).Value = 4
    $ws.Cells.Item(25, 8).Value = 10
    $ws.Cells.Item(25, 9).Value = 'boolean'
    $ws.Cells.Item(25, 10).Value = 'workflow'
    $ws.Cells.Item(25, 12).Value = 0
    $ws.Cells.Item(25, 14).Value = '[{"value":"1","label":"مراجعه کرد"},{"value":"0","label":"مراجعه نکرد"}]'
    # Row 26: VAR_UID=9244949946521480773d796066448722
    $ws.Cells.Item(26, 1).Value = 26
    $ws.Cells.Item(26, 2).Value = '9244949946521480773d796066448722'
    $ws.Cells.Item(26, 3).Value = '61815949564e4999e9835d8060262100'
    $ws.Cells.Item(26, 4).Value = 2
    $ws.Cells.Item(26, 5).Value = 'referral_and_notice_sms'
    $ws.Cells.Item(26, 6).Value = 'string'
    $ws.Cells.Item(26, 7).Value = 1
    $ws.Cells.Item(26, 8).Value = 10
    $ws.Cells.Item(26, 9).Value = 'string'
    $ws.Cells.Item(26, 10).Value = 'workflow'
    $ws.Cells.Item(26, 12).Value = 0
    $ws.Cells.Item(26, 14).Value = '[]'
    # Row 27: VAR_UID=94172160065211a796e7149060650444
    $ws.Cells.Item(27, 1).Value = 23
    $ws.Cells.Item(27, 2).Value = '94172160065211a796e7149060650444'
    $ws.Cells.Item(27, 3).Value = '61815949564e4999e9835d8060262100'
    $ws.Cells.Item(27, 4).Value = 2
    $ws.Cells.Item(27, 5).Value = 'warning_deadline'
    $ws.Cells.Item(27, 6).Value = 'string'
    $ws.Cells.Item(27, 7).Value = 1
    $ws.Cells.Item(27, 8).Value = 10
    $ws.Cells.Item(27, 9).Value = 'string'
    $ws.Cells.Item(27, 10).Value = 'workflow'
    $ws.Cells.Item(27, 12).Value = 0
    $ws.Cells.Item(27, 14).Value = '[]'
    # Row 28: VAR_UID=97087285364e34e37455df0099041822
    $ws.Cells.Item(28, 1).Value = 8
    $ws.Cells.Item(28, 2).Value = '97087285364e34e37455df0099041822'
    $ws.Cells.Item(28, 3).Value = '20109551764e348a7a8c913045934777'
    $ws.Cells.Item(28, 4).Value = 1
    $ws.Cells.Item(28, 5).Value = 'description'
    $ws.Cells.Item(28, 6).Value = 'string'
    $ws.Cells.Item(28, 7).Value = 1
    $ws.Cells.Item(28, 8).Value = 10
    $ws.Cells.Item(28, 9).Value = 'string'
    $ws.Cells.Item(28, 10).Value = 'workflow'
    $ws.Cells.Item(28, 12).Value = 0
    $ws.Cells.Item(28, 14).Value = '[]'
    # Row 29: VAR_UID=97160399964e35567ca3989086381213
    $ws.Cells.Item(29, 1).Value = 9
    $ws.Cells.Item(29, 2).Value = '97160399964e35567ca3989086381213'
    $ws.Cells.Item(29, 3).Value = '20109551764e348a7a8c913045934777'
    $ws.Cells.Item(29, 4).Value = 1
    $ws.Cells.Item(29, 5).Value = 'name'
    $ws.Cells.Item(29, 6).Value = 'string'
    $ws.Cells.Item(29, 7).Value = 1
    $ws.Cells.Item(29, 8).Value = 10
    $ws.Cells.Item(29, 9).Value = 'string'
    $ws.Cells.Item(29, 10).Value = 'workflow'
    $ws.Cells.Item(29, 12).Value = 0
    $ws.Cells.Item(29, 14).Value = '[]'
    # Row 30: VAR_UID=995623521652147b8e3e140020518509
    $ws.Cells.Item(30, 1).Value = 25
    $ws.Cells.Item(30, 2).Value = '995623521652147b8e3e140020518509'
    $ws.Cells.Item(30, 3).Value = '61815949564e4999e9835d8060262100'
    $ws.Cells.Item(30, 4).Value = 2
    $ws.Cells.Item(30, 5).Value = 'inactivity_commitment'
    $ws.Cells.Item(30, 6).Value = 'boolean'
    $ws.Cells.Item(30, 7).Value = 4
    $ws.Cells.Item(30, 8).Value = 10
    $ws.Cells.Item(30, 9).Value = 'boolean'
    $ws.Cells.Item(30, 10).Value = 'workflow'
    $ws.Cells.Item(30, 12).Value = 0
    $ws.Cells.Item(30, 14).Value = '[{"value":"1","label":"دارد"},{"value":"0","label":"ندارد"}]'
